$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2:B26").Value = [double]"0.9999824846980448"
$ws.Range("C2:C26").Value = [double]"0.9990072066285755"
$ws.Range("D2:D26").Value = [double]"0.9999880092558681"
$ws.Range("E2:E26").Value = [double]"0.9999998552291943"
$ws.Range("F2:F26").Value = [double]"0.9999921855090684"
$ws.Range("G2:G26").Value = [double]"1.634976554226276e-05"
$ws.Range("H2:H26").Value = [double]"0.0009267290336323482"
$ws.Range("I2:I26").Value = [double]"1.445341323746581e-05"
$ws.Range("J2:J26").Value = [double]"7.975904722099923e-08"
$ws.Range("K2:K26").Value = [double]"7.266586142343406e-06"
$ws.Range("L2:L26").Value = [double]"0.0002547310818618769"
$ws.Range("M2:M26").Value = [double]"0.004043484331892824"
$ws.Range("N2:N26").Value = [double]"1.000016814689877"
$ws.Range("O2:O26").Value = [double]"0.004215623831067973"
$ws.Range("P2:P26").Value = [double]"120.0425940012932"
$ws.Range("Q2:Q26").Value = [double]"179.767509419835"
